# Set Property sheet's "Private" column (D) to TRUE for every data row,
# per commit message "set property's private value as true".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# D2:D6 holds the "Private" boolean flag for each row -> flip to TRUE
$privateRange = $ws.Range("D2:D6")
$privateRange.Value = $true

# Re-apply the TRUE/FALSE list validation over the full column F range
# (consolidating it into one contiguous rule) and extend the same rule to
# the newly-edited column D range.
$fRange = $ws.Range("F2:F1048576")
$fRange.Validation.Delete()
$fRange.Validation.Add(3, 1, 1, "TRUE,FALSE")
$privateRange.Validation.Add(3, 1, 1, "TRUE,FALSE")

# Reflect where the edit was made in the sheet's active selection.
[void]$privateRange.Select()

$wb.Save()
